$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M), Std_Experiment_Time (N)
# which were computed incorrectly (off-by-one / swapped averaging).
$ws.Range("G2").Value = 7.40809352
$ws.Range("H2").Value = 398.14300002
$ws.Range("M2").Value = 1.010627121043755
$ws.Range("N2").Value = 105.4411928558437
$ws.Range("G3").Value = 11.54658394
$ws.Range("H3").Value = 1036.00402306
$ws.Range("M3").Value = 1.594986336144876
$ws.Range("N3").Value = 293.9491129333546
$ws.Range("G4").Value = 2.43512314
$ws.Range("H4").Value = 72.60676776
$ws.Range("M4").Value = 0.4161116558957283
$ws.Range("N4").Value = 25.90025464032825
$ws.Range("G5").Value = 3.2722633
$ws.Range("H5").Value = 163.49937598
$ws.Range("M5").Value = 0.562781947459656
$ws.Range("N5").Value = 60.38638085824623
$ws.Range("G6").Value = 0.7669273600000001
$ws.Range("H6").Value = 12.3088315
$ws.Range("M6").Value = 0.1934415784073938
$ws.Range("N6").Value = 6.745546862543131
$ws.Range("G7").Value = 0.92498132
$ws.Range("H7").Value = 24.45588658
$ws.Range("M7").Value = 0.171736857617654
$ws.Range("N7").Value = 9.945613175554639
$ws.Range("G8").Value = 0.40577266
$ws.Range("H8").Value = 4.36777126
$ws.Range("M8").Value = 0.09406143466749706
$ws.Range("N8").Value = 2.226559821479375
$ws.Range("G9").Value = 0.43350224
$ws.Range("H9").Value = 7.997219659999999
$ws.Range("M9").Value = 0.09329857909907424
$ws.Range("N9").Value = 3.844196882024126
$ws.Range("G10").Value = 0.2449111
$ws.Range("H10").Value = 1.95037218
$ws.Range("M10").Value = 0.05633855779818613
$ws.Range("N10").Value = 0.939673584010758
$ws.Range("G11").Value = 0.24452178
$ws.Range("H11").Value = 3.54891342
$ws.Range("M11").Value = 0.06066215385272521
$ws.Range("N11").Value = 2.040596071443755
$ws.Range("G12").Value = 0.17271368
$ws.Range("H12").Value = 1.11245764
$ws.Range("M12").Value = 0.04531058831315445
$ws.Range("N12").Value = 0.6068169002690471
$ws.Range("G13").Value = 0.15489318
$ws.Range("H13").Value = 1.81991824
$ws.Range("M13").Value = 0.0408042252282844
$ws.Range("N13").Value = 1.088608759030616
